$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.999999998154007196937698154215
$ws.Range("C2").Value = 0.000000000319760100408631121627
$ws.Range("D2").Value = 0.000000000282360682425450125399
$ws.Range("E2").Value = 0.000000001243880157806520025340
$ws.Range("F2").Value = 46052
